$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: update the base lat/lon input values ---
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = 45
$ws.Range("C3").Value = -85
$ws.Range("D3").Value = -55

# --- Row 8: J8 formula changed from "=4*12" to "=B5*4" ---
$ws.Range("J8").Formula = "=B5*4"

# --- Rows 9-14: updated refinement-level counts (B column / J column) ---
$ws.Range("B9").Value = 2
$ws.Range("J9").Value = 2

$ws.Range("B10").Value = 2
$ws.Range("J10").Value = 2

$ws.Range("B12").Value = 11
$ws.Range("J12").Value = 11

$ws.Range("B13").Value = 4
$ws.Range("J13").Value = 4

$ws.Range("B14").Value = 6
$ws.Range("J14").Value = 6

# --- New rows 28-30: extra lat offset calculations ---
$ws.Range("E28").Formula = "=300/3600/60"
$ws.Range("G28").Value = " 40.819273°"

$ws.Range("D29").Formula = "=-72.56"
$ws.Range("E29").Formula = "=D29+E28"
$ws.Range("E29").NumberFormat = "0.000000000"
$ws.Range("G29").Value = " 40.817693°"

$ws.Range("E30").Formula = "=D29-E28"
$ws.Range("E30").NumberFormat = "0.000000000"

# --- Column width for column E + selection to match the final view state ---
$ws.Columns.Item(5).ColumnWidth = 13.1640625
$ws.Range("E29").Select()
